$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")
$ws.Activate()

$ws.Range("C32").Value = 0.8046875
$ws.Range("G32").Value = 42541
$ws.Range("H32").Value = "13.12.txt"
$ws.Range("G33").Value = 42541

$ws.Range("C32").Select()
